$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 388.36
$ws.Range("I28").Value = 347.4737
$ws.Range("J28").Value = 517.8333
$ws.Range("K28").Value = 347.4737
$ws.Range("L28").Value = 517.8333
$ws.Range("M28").Value = 137.5263
$ws.Range("N28").Value = -1487.8333

$ws.Range("H29").Value = 2106.389
$ws.Range("J29").Value = 3512.6667
$ws.Range("L29").Value = 10538.0001
$ws.Range("N29").Value = -11100.0001

$ws.Range("H33").Value = 96.85714
$ws.Range("I33").Value = 89.666664
$ws.Range("K33").Value = 89.666664
$ws.Range("M33").Value = 139.333336

$ws.Range("H43").Value = 9749.5
$ws.Range("I43").Value = 12666
$ws.Range("K43").Value = 12666
$ws.Range("M43").Value = -12597

$ws.Range("H80").Value = 5450.75
$ws.Range("I80").Value = 8291.200000000001
$ws.Range("J80").Value = 716.6667
$ws.Range("K80").Value = 24873.6
$ws.Range("L80").Value = 2150.0001
$ws.Range("M80").Value = -23875.6
$ws.Range("N80").Value = -4146.0001

$ws.Range("H83").Value = 5450.75
$ws.Range("I83").Value = 8291.200000000001
$ws.Range("J83").Value = 716.6667
$ws.Range("K83").Value = 74620.8
$ws.Range("L83").Value = 6450.0003
$ws.Range("M83").Value = -69628.8
$ws.Range("N83").Value = -16434.0003

$ws.Range("H125").Value = 7924.7144
$ws.Range("I125").Value = 7271.6665
$ws.Range("J125").Value = 8414.5
$ws.Range("K125").Value = 65444.9985
$ws.Range("L125").Value = 75730.5
$ws.Range("M125").Value = -62984.9985
$ws.Range("N125").Value = -80650.5

$ws.Range("H132").Value = 2264.5908
$ws.Range("I132").Value = 1201.1578
$ws.Range("K132").Value = 3603.4734
$ws.Range("M132").Value = -1073.4734

$ws.Range("H135").Value = 1682.7
$ws.Range("I135").Value = 1704.5
$ws.Range("K135").Value = 15340.5
$ws.Range("M135").Value = -12805.5

$ws.Range("H137").Value = 7965.6665
$ws.Range("I137").Value = 8711.375
$ws.Range("K137").Value = 26134.125
$ws.Range("M137").Value = -23584.125

$ws.Range("H138").Value = 3314.7097
$ws.Range("J138").Value = 2893.158
$ws.Range("L138").Value = 8679.474
$ws.Range("N138").Value = -18959.474

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 258708.95
$ws.Range("I32").Value = 992.4545000000001
$ws.Range("J32").Value = 1676149.6
$ws.Range("K32").Value = 992.4545000000001
$ws.Range("L32").Value = 1676149.6
$ws.Range("M32").Value = -705.4545000000001
$ws.Range("N32").Value = -1676723.6

$ws.Range("H61").Value = 6085.7646
$ws.Range("I61").Value = 6130.533
$ws.Range("K61").Value = 6130.533
$ws.Range("M61").Value = -5918.533

$ws.Range("H110").Value = 1449.4231
$ws.Range("I110").Value = 723.9167
$ws.Range("J110").Value = 2071.2856
$ws.Range("K110").Value = 723.9167
$ws.Range("L110").Value = 2071.2856
$ws.Range("M110").Value = 1321.0833
$ws.Range("N110").Value = -6161.2856

$ws.Range("H123").Value = 65499.75
$ws.Range("J123").Value = 65499.75
$ws.Range("L123").Value = 65499.75
$ws.Range("N123").Value = -75299.75

$ws.Range("H132").Value = 3278.2856
$ws.Range("I132").Value = 2585.7144
$ws.Range("K132").Value = 7757.1432
$ws.Range("M132").Value = -5227.1432

$ws.Range("H136").Value = 6085.7646
$ws.Range("I136").Value = 6130.533
$ws.Range("K136").Value = 18391.599
$ws.Range("M136").Value = -15841.599

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13022.818
$ws.Range("I31").Value = 17783.834
$ws.Range("J31").Value = 7309.6
$ws.Range("K31").Value = 17783.834
$ws.Range("L31").Value = 7309.6
$ws.Range("M31").Value = -17488.834
$ws.Range("N31").Value = -7899.6

$ws.Range("H34").Value = 13022.818
$ws.Range("I34").Value = 17783.834
$ws.Range("J34").Value = 7309.6
$ws.Range("K34").Value = 17783.834
$ws.Range("L34").Value = 7309.6
$ws.Range("M34").Value = -17581.834
$ws.Range("N34").Value = -7713.6

$ws.Range("H51").Value = 22233
$ws.Range("J51").Value = 22979.6
$ws.Range("L51").Value = 22979.6
$ws.Range("N51").Value = -24451.6

$ws.Range("H58").Value = 2409.2593
$ws.Range("I58").Value = 1829.5264
$ws.Range("K58").Value = 1829.5264
$ws.Range("M58").Value = -1626.5264

$ws.Range("H61").Value = 22233
$ws.Range("J61").Value = 22979.6
$ws.Range("L61").Value = 22979.6
$ws.Range("N61").Value = -23675.6

$ws.Range("H68").Value = 63333
$ws.Range("J68").Value = 63333
$ws.Range("L68").Value = 63333
$ws.Range("N68").Value = -64831

$ws.Range("H71").Value = 63333
$ws.Range("J71").Value = 63333
$ws.Range("L71").Value = 189999
$ws.Range("N71").Value = -197487

$ws.Range("H74").Value = 30057
$ws.Range("J74").Value = 30057
$ws.Range("L74").Value = 30057
$ws.Range("N74").Value = -31805

$ws.Range("H77").Value = 30057
$ws.Range("J77").Value = 30057
$ws.Range("L77").Value = 90171
$ws.Range("N77").Value = -98907

$ws.Range("H107").Value = 1166.8823
$ws.Range("I107").Value = 790.7
$ws.Range("J107").Value = 1323.625
$ws.Range("K107").Value = 790.7
$ws.Range("L107").Value = 1323.625
$ws.Range("M107").Value = 1129.3
$ws.Range("N107").Value = -5163.625

$ws.Range("H132").Value = 5327.706
$ws.Range("I132").Value = 3596.625
$ws.Range("K132").Value = 10789.875
$ws.Range("M132").Value = -8259.875

$ws.Range("H134").Value = 1440.2941
$ws.Range("I134").Value = 1159.8572
$ws.Range("K134").Value = 3479.5716
$ws.Range("M134").Value = -944.5715999999998

$ws.Range("H136").Value = 2409.2593
$ws.Range("I136").Value = 1829.5264
$ws.Range("K136").Value = 5488.5792
$ws.Range("M136").Value = -2938.5792

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 272.30768
$ws.Range("I2").Value = 29.428572
$ws.Range("J2").Value = 555.6667
$ws.Range("K2").Value = 176.571432
$ws.Range("L2").Value = 3334.0002
$ws.Range("M2").Value = -63.57143199999999
$ws.Range("N2").Value = -3560.0002

$ws.Range("H7").Value = 5882470
$ws.Range("I7").Value = 8333440
$ws.Range("J7").Value = 142.4
$ws.Range("K7").Value = 25000320
$ws.Range("L7").Value = 427.2
$ws.Range("M7").Value = -25000208
$ws.Range("N7").Value = -651.2

$ws.Range("H12").Value = 216.58824
$ws.Range("I12").Value = 250.27272
$ws.Range("J12").Value = 154.83333
$ws.Range("K12").Value = 750.81816
$ws.Range("L12").Value = 464.49999
$ws.Range("M12").Value = -577.81816
$ws.Range("N12").Value = -810.49999

$ws.Range("H23").Value = 113033.11
$ws.Range("I23").Value = 2599.5
$ws.Range("J23").Value = 144585.58
$ws.Range("K23").Value = 7798.5
$ws.Range("L23").Value = 433756.74
$ws.Range("M23").Value = -7563.5
$ws.Range("N23").Value = -434226.74

$ws.Range("H34").Value = 1175.2667
$ws.Range("I34").Value = 792.125
$ws.Range("J34").Value = 2707.8333
$ws.Range("K34").Value = 2376.375
$ws.Range("L34").Value = 8123.499899999999
$ws.Range("M34").Value = -2292.375
$ws.Range("N34").Value = -8291.499899999999

$ws.Range("H39").Value = 3400
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 3400
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 10200
$ws.Range("N39").Value = -10788
$ws.Range("M39").ClearContents()

$ws.Range("H55").Value = 4339.8
$ws.Range("J55").Value = 4237.25
$ws.Range("L55").Value = 12711.75
$ws.Range("N55").Value = -13065.75

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1550
$ws.Range("I107").Value = 1550
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1550
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 370
$ws.Range("N107").ClearContents()

$ws.Range("H132").Value = 4169
$ws.Range("I132").Value = 2750
$ws.Range("J132").Value = 4878.5
$ws.Range("K132").Value = 8250
$ws.Range("L132").Value = 14635.5
$ws.Range("M132").Value = -5720
$ws.Range("N132").Value = -19695.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5963.5
$ws.Range("I7").Value = 5957.8
$ws.Range("K7").Value = 5957.8
$ws.Range("M7").Value = -5845.8

$ws.Range("H126").Value = 5963.5
$ws.Range("I126").Value = 5957.8
$ws.Range("K126").Value = 17873.4
$ws.Range("M126").Value = -15403.4
